$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A167").Value = "2023-12-10 17:12:59"
$ws.Range("B167").Value = 0.0006000000000000001

$ws.Range("A168").Value = "2023-12-10 17:13:10"
$ws.Range("B168").Value = 0.0006000000000000001
